$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('clinics', 'さくら歯科'),
    @('clinics', 'たんぽぽ歯科'),
    @('clinics', 'ありす歯科'),
    @('clinics', '春日井きらり歯科'),
    @('clinics', '松戸ありす歯科'),
    @('clinics', '池下さくら歯科'),
    @('clinics', '日進赤池たんぽぽ歯科'),
    @('clinics', '春日井アップル歯科'),
    @('clinics', '緑区さくら医院・歯科'),
    @('clinics', '緑区さくら医院・医科'),
    @('clinics', '緑区さくら医院・健診'),
    @('clinics', '金沢さくら医院・健診'),
    @('clinics', '金沢さくら医院・婦人科'),
    @('clinics', 'ハピネス歯科'),
    @('clinics', '流山ハピネス歯科'),
    @('clinics', 'イーアス春日井'),
    @('clinics', '名駅さくら医院・名古屋歯科・歯科'),
    @('clinics', '名駅さくら医院・名古屋歯科・内科'),
    @('clinics', '名駅さくら医院・名古屋歯科・皮膚科'),
    @('clinics', 'きらり大森歯科'),
    @('clinics', 'クローバー歯科'),
    @('clinics', '流山ありす歯科・矯正歯科'),
    @('clinics', '長久手さくら歯科・矯正歯科'),
    @('clinics', 'ヒロデンタル'),
    @('clinics', '町屋さくら歯科・矯正歯科'),
    @('clinics', '金沢さくら医院 歯科'),
    @('clinics', 'きた矯正歯科クリニック'),
    @('clinics', 'デンタルオフィス増田'),
    @('clinics', '祥南歯科・矯正歯科医院'),
    @('clinics', '三田矯正歯科医院'),
    @('clinics', 'グランド歯科医院'),
    @('clinics', '訪問部（春日井事務所）'),
    @('clinics', '訪問部（杁中事務所）'),
    @('clinics', '春日井事務局'),
    @('clinics', '名古屋駅JPタワー事務所'),
    @('clinics', 'いりなか事務所'),
    @('clinics', '池下事務所'),
    @('clinics', '流山事務所'),
    @('clinics', '浄心事務所'),
    @('positions', '歯科医師'),
    @('positions', '歯科衛生士'),
    @('positions', '歯科助手'),
    @('positions', '受付'),
    @('positions', '事務'),
    @('positions', 'マネージャー'),
    @('status', 'active'),
    @('status', 'leave'),
    @('status', 'retired'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

Write-Output "rows written"
